# UCD_EngArch_Path_Civil_ME_Modules.xlsx - "updated slide templates + paths"
#
# The real-world commit corrected the Stage column (E) for the Civil/ME
# module rows, which had mistakenly been filled with the literal text "M"
# (a stray shared string) instead of the numeric Stage. The Major column
# (F), which legitimately holds the "T298" programme code, is left
# untouched content-wise; once the stray "M" shared string is no longer
# referenced anywhere it drops out of the shared-string table and every
# subsequent string index shifts down by one automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage (column E) should be numeric 4 for rows 38-44 ...
$ws.Range("E38:E44").Value = 4

# ... and numeric 5 for rows 45-52.
$ws.Range("E45:E52").Value = 5

# Restore the view state (scroll position + selection) to where the author
# left off editing this block of rows.
$ws.Range("E45").Select()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
